$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.790.70"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "3.554.03"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.85"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D7").Value = "3.552.24"
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +5.88%  "
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.81"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "4.160.24"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000197"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.09"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "3.560.13"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.791.56"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.117"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.01"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.26"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.75"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "427.27"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.600"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.85"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").Value = "3.697.69"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -4.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.05"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "3.563.47"
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.51"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("E36").Value = "  -8.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.66"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "177.36"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.31"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0830"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.04"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.99%  "
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.53"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("E45").Value = "  -6.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("E47").Value = "  -6.01%  "
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.44"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("E50").Value = "  -4.66%  "
$ws.Range("E51").Value = "  -3.18%  "
